# new-hampshire_overview.xlsx edit script
# - Inserts a new leading data column "Share of 990 filers with government
#   grants at risk" (previously the last column) in every sheet, shifting the
#   remaining metric columns one slot to the right.
# - Renames "Operating surplus with/without government grants (%)" to
#   "Size of operating surplus with/without government grants".
# - Relabels a handful of row headers (Congressional District, Size buckets,
#   Subsector names) and re-sorts rows alphabetically by that label where the
#   source data changed order.
#
# Helper: write a value into a cell while forcing it to be stored as literal
# text (matches the source workbook, where every value - including things
# that look like numbers, currency or percentages - is inline text, not a
# real number/currency/percentage). We flip the number format to Text before
# the write so the engine doesn't "smart type" the string, then reset the
# cell style back to Normal so no stray formatting sticks around.
function Set-TextCell {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Helper: write a header-row value (row 1). These are plain words, so no
# special text-forcing is required, and the existing bold/centered header
# style (s="1") is left untouched automatically since we only touch .Value.
function Set-HeaderCell {
    param($ws, $row, $col, $val)
    $ws.Cells.Item($row, $col).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overall
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overall")

Set-HeaderCell $ws 1 1 "Share of 990 filers with government grants at risk"
Set-HeaderCell $ws 1 2 "Number of 990 filers with government grants"
Set-HeaderCell $ws 1 3 "Total government grants (`$)"
Set-HeaderCell $ws 1 4 "Size of operating surplus with government grants"
Set-HeaderCell $ws 1 5 "Size of operating surplus without government grants"

Set-TextCell $ws 2 1 "62.90%"
Set-TextCell $ws 2 2 "663"
Set-TextCell $ws 2 3 "`$1,133,383,584"
Set-TextCell $ws 2 4 "10.55%"
Set-TextCell $ws 2 5 "-7.39%"

# ---------------------------------------------------------------------
# Sheet 2: County
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("County")

Set-HeaderCell $ws 1 1 "Geography"
Set-HeaderCell $ws 1 2 "Share of 990 filers with government grants at risk"
Set-HeaderCell $ws 1 3 "Number of 990 filers with government grants"
Set-HeaderCell $ws 1 4 "Total government grants (`$)"
Set-HeaderCell $ws 1 5 "Size of operating surplus with government grants"
Set-HeaderCell $ws 1 6 "Size of operating surplus without government grants"

$countyRows = @(
    @("United States",      "67.35%", "103,475", "`$267,700,640,005", "9.05%",  "-12.83%"),
    @("New Hampshire",      "62.90%", "663",     "`$1,133,383,584",   "10.55%", "-7.39%"),
    @("Belknap County",     "44.83%", "29",      "`$11,058,626",      "12.78%", "0.76%"),
    @("Carroll County",     "56.36%", "55",      "`$19,490,173",      "15.35%", "-2.36%"),
    @("Cheshire County",    "69.23%", "39",      "`$45,943,451",      "12.71%", "-7.58%"),
    @("Coos County",        "66.67%", "18",      "`$46,524,474",      "8.52%",  "-27.19%"),
    @("Grafton County",     "58.67%", "75",      "`$298,933,446",     "10.48%", "-1.82%"),
    @("Hillsborough County","66.67%", "162",     "`$407,439,385",     "9.25%",  "-10.30%"),
    @("Merrimack County",   "67.80%", "118",     "`$165,086,443",     "6.37%",  "-13.16%"),
    @("Rockingham County",  "59.41%", "101",     "`$62,471,837",      "13.93%", "-4.80%"),
    @("Strafford County",   "64.58%", "48",      "`$67,164,393",      "14.72%", "-6.79%"),
    @("Sullivan County",    "61.11%", "18",      "`$9,271,356",       "14.71%", "-1.91%")
)

$r = 2
foreach ($row in $countyRows) {
    Set-TextCell $ws $r 1 $row[0]
    Set-TextCell $ws $r 2 $row[1]
    Set-TextCell $ws $r 3 $row[2]
    Set-TextCell $ws $r 4 $row[3]
    Set-TextCell $ws $r 5 $row[4]
    Set-TextCell $ws $r 6 $row[5]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 3: Congressional District
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Congressional District")

Set-HeaderCell $ws 1 1 "Geography"
Set-HeaderCell $ws 1 2 "Share of 990 filers with government grants at risk"
Set-HeaderCell $ws 1 3 "Number of 990 filers with government grants"
Set-HeaderCell $ws 1 4 "Total government grants (`$)"
Set-HeaderCell $ws 1 5 "Size of operating surplus with government grants"
Set-HeaderCell $ws 1 6 "Size of operating surplus without government grants"

$cdRows = @(
    @("United States",           "67.35%", "103,475", "`$267,700,640,005", "9.05%",  "-12.83%"),
    @("New Hampshire",           "62.90%", "663",     "`$1,133,383,584",   "10.55%", "-7.39%"),
    @("Congressional District 1","63.64%", "319",     "`$495,923,217",     "10.76%", "-7.71%"),
    @("Congressional District 2","62.21%", "344",     "`$637,460,367",     "10.39%", "-7.00%")
)

$r = 2
foreach ($row in $cdRows) {
    Set-TextCell $ws $r 1 $row[0]
    Set-TextCell $ws $r 2 $row[1]
    Set-TextCell $ws $r 3 $row[2]
    Set-TextCell $ws $r 4 $row[3]
    Set-TextCell $ws $r 5 $row[4]
    Set-TextCell $ws $r 6 $row[5]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 4: Size
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Size")

Set-HeaderCell $ws 1 1 "Size"
Set-HeaderCell $ws 1 2 "Share of 990 filers with government grants at risk"
Set-HeaderCell $ws 1 3 "Number of 990 filers with government grants"
Set-HeaderCell $ws 1 4 "Total government grants (`$)"
Set-HeaderCell $ws 1 5 "Size of operating surplus with government grants"
Set-HeaderCell $ws 1 6 "Size of operating surplus without government grants"

$sizeRows = @(
    @("Between `$100K and `$499K", "63.00%", "227", "`$23,135,877",  "11.57%", "-11.49%"),
    @("Between `$1M and `$4.99M",  "61.93%", "176", "`$130,239,318", "13.24%", "-6.72%"),
    @("Between `$500K and `$999K", "67.26%", "113", "`$31,835,568",  "10.96%", "-12.41%"),
    @("Between `$5M and `$9.99M",  "59.26%", "27",  "`$50,104,415",  "9.97%",  "-8.91%"),
    @("Greater than `$10M",        "65.26%", "95",  "`$896,868,088", "5.26%",  "-3.00%"),
    @("Less than `$100K",          "44.00%", "25",  "`$1,200,318",   "28.69%", "2.38%"),
    @("Total",                    "62.90%", "663", "`$1,133,383,584","10.55%", "-7.39%")
)

$r = 2
foreach ($row in $sizeRows) {
    Set-TextCell $ws $r 1 $row[0]
    Set-TextCell $ws $r 2 $row[1]
    Set-TextCell $ws $r 3 $row[2]
    Set-TextCell $ws $r 4 $row[3]
    Set-TextCell $ws $r 5 $row[4]
    Set-TextCell $ws $r 6 $row[5]
    $r++
}

# ---------------------------------------------------------------------
# Sheet 5: Subsector
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Subsector")

Set-HeaderCell $ws 1 1 "Subsector"
Set-HeaderCell $ws 1 2 "Share of 990 filers with government grants at risk"
Set-HeaderCell $ws 1 3 "Number of 990 filers with government grants"
Set-HeaderCell $ws 1 4 "Total government grants (`$)"
Set-HeaderCell $ws 1 5 "Size of operating surplus with government grants"
Set-HeaderCell $ws 1 6 "Size of operating surplus without government grants"

$subsectorRows = @(
    @("Arts, Culture, and Humanities",   "51.85%", "54",  "`$10,875,853",  "23.59%", "-0.26%"),
    @("Education (Excluding Universities)","63.95%","86", "`$60,703,868",  "10.44%", "-13.09%"),
    @("Environment and Animals",         "33.33%", "36",  "`$5,010,754",   "30.15%", "20.72%"),
    @("Health (Excluding Hospitals)",    "76.56%", "64",  "`$102,493,282", "10.26%", "-11.81%"),
    @("Hospitals",                       "75.00%", "4",   "`$1,176,777",   "-5.63%", "-12.68%"),
    @("Human Services",                  "67.30%", "211", "`$157,151,757", "7.37%",  "-11.10%"),
    @("International, Foreign Affairs",  "33.33%", "3",   "`$171,352",     "15.31%", "7.12%"),
    @("Public, Societal Benefit",        "55.00%", "40",  "`$28,513,999",  "18.30%", "-0.30%"),
    @("Religion Related",                "33.33%", "6",   "`$453,109",     "29.21%", "11.80%"),
    @("Unclassified",                    "65.33%", "150", "`$494,942,147", "10.40%", "-8.72%"),
    @("Universities",                    "55.56%", "9",   "`$271,890,686", "9.63%",  "-1.15%"),
    @("Total",                           "62.90%", "663", "`$1,133,383,584","10.55%","-7.39%")
)

$r = 2
foreach ($row in $subsectorRows) {
    Set-TextCell $ws $r 1 $row[0]
    Set-TextCell $ws $r 2 $row[1]
    Set-TextCell $ws $r 3 $row[2]
    Set-TextCell $ws $r 4 $row[3]
    Set-TextCell $ws $r 5 $row[4]
    Set-TextCell $ws $r 6 $row[5]
    $r++
}
